# rename backend to en-us
# Renames the sheets and translates header labels from Portuguese to
# English, and strips the header styling (bold/fill/border/alignment),
# the frozen header row pane, and the custom column widths that were
# tied to that styling.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (Portuguese -> English) ---------------------------------
$sheetRenames = @{
    "Usuarios"      = "Users"
    "Produtos"      = "Products"
    "Movimentacoes" = "Movements"
    "Vendas"        = "Sales"
}

foreach ($ws in $wb.Worksheets) {
    $newName = $sheetRenames[$ws.Name]
    if ($newName) {
        $ws.Name = $newName
    }
}

# --- Header translations per sheet -----------------------------------------
$headerTranslations = @{
    "Users" = @("email", "role", "createdAt")
    "Products" = @("id", "name", "description", "price", "stock", "imageId", "active", "createdAt", "updatedAt")
    "Movements" = @("id", "productId", "type", "quantity", "reason", "userEmail", "createdAt")
    "Sales" = @("id", "productId", "quantity", "unitPrice", "total", "paymentMethod", "userEmail", "createdAt")
}

foreach ($ws in $wb.Worksheets) {
    $headers = $headerTranslations[$ws.Name]
    if ($headers) {
        for ($i = 0; $i -lt $headers.Length; $i++) {
            $cell = $ws.Cells.Item(1, $i + 1)
            $cell.Value = $headers[$i]
            # Strip the bold/fill/border/center-align header style back to
            # the workbook default.
            $cell.ClearFormats()
        }

        # Reset the columns that previously carried an explicit custom
        # width tied to the header styling.
        for ($i = 0; $i -lt $headers.Length; $i++) {
            $ws.Columns.Item($i + 1).ClearFormats()
        }

        # Unfreeze the header row (remove the frozen pane / simplify the
        # selection back to a single top-left cell).
        $ws.Activate()
        $excel.ActiveWindow.FreezePanes = $false
    }
}

# Restore the originally active sheet/selection (first tab).
$wb.Worksheets.Item(1).Activate()
